$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "63.496.88"
$ws.Cells.Item(2,5).Value = "  -2.23%  "

$ws.Cells.Item(3,4).Value = "3.477.15"
$ws.Cells.Item(3,5).Value = "  -1.13%  "

$ws.Cells.Item(4,5).Value = "  -0.12%  "

$c = $ws.Cells.Item(5,4)
$c.NumberFormat = "@"
$c.Value = "581.25"
$c.ClearFormats()
$ws.Cells.Item(5,5).Value = "  -2.68%  "

$c = $ws.Cells.Item(6,4)
$c.NumberFormat = "@"
$c.Value = "131.08"
$c.ClearFormats()
$ws.Cells.Item(6,5).Value = "  -1.97%  "

$ws.Cells.Item(7,4).Value = "3.476.01"
$ws.Cells.Item(7,5).Value = "  -1.13%  "

$ws.Cells.Item(8,5).Value = "  -0.03%  "

$ws.Cells.Item(9,5).Value = "  -2.02%  "

$ws.Cells.Item(10,5).Value = "  -0.44%  "

$c = $ws.Cells.Item(11,4)
$c.NumberFormat = "@"
$c.Value = "7.13"
$c.ClearFormats()
$ws.Cells.Item(11,5).Value = "  +0.07%  "

$c = $ws.Cells.Item(12,4)
$c.NumberFormat = "@"
$c.Value = "0.379"
$c.ClearFormats()
$ws.Cells.Item(12,5).Value = "  -0.70%  "

$ws.Cells.Item(13,4).Value = "4.076.74"
$ws.Cells.Item(13,5).Value = "  -1.24%  "

$c = $ws.Cells.Item(14,4)
$c.NumberFormat = "@"
$c.Value = "27.35"
$c.ClearFormats()
$ws.Cells.Item(14,5).Value = "  +0.18%  "

$ws.Cells.Item(15,5).Value = "  +1.60%  "

$ws.Cells.Item(16,4).Value = "3.504.18"
$ws.Cells.Item(16,5).Value = "  -0.64%  "

$ws.Cells.Item(17,5).Value = "  -2.81%  "

$ws.Cells.Item(18,4).Value = "63.626.87"
$ws.Cells.Item(18,5).Value = "  -2.05%  "

$c = $ws.Cells.Item(19,4)
$c.NumberFormat = "@"
$c.Value = "9.98"
$c.ClearFormats()
$ws.Cells.Item(19,5).Value = "  -0.13%  "

$ws.Cells.Item(20,5).Value = "  -0.56%  "

$ws.Cells.Item(21,5).Value = "  -0.26%  "

$c = $ws.Cells.Item(22,4)
$c.NumberFormat = "@"
$c.Value = "382.64"
$c.ClearFormats()
$ws.Cells.Item(22,5).Value = "  -2.04%  "

$ws.Cells.Item(23,5).Value = "  +0.66%  "

$ws.Cells.Item(24,4).Value = "3.617.78"
$ws.Cells.Item(24,5).Value = "  -1.21%  "

$c = $ws.Cells.Item(25,4)
$c.NumberFormat = "@"
$c.Value = "72.83"
$c.ClearFormats()
$ws.Cells.Item(25,5).Value = "  -1.58%  "

$c = $ws.Cells.Item(26,4)
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.ClearFormats()
$ws.Cells.Item(26,5).Value = "  -0.12%  "

$c = $ws.Cells.Item(27,4)
$c.NumberFormat = "@"
$c.Value = "0.0000113"
$c.ClearFormats()
$ws.Cells.Item(27,5).Value = "  -0.34%  "

$c = $ws.Cells.Item(28,4)
$c.NumberFormat = "@"
$c.Value = "1.56"
$c.ClearFormats()
$ws.Cells.Item(28,5).Value = "  -4.00%  "

$ws.Cells.Item(29,2).Value = "RenderToken"
$ws.Cells.Item(29,3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Cells.Item(29,4)
$c.NumberFormat = "@"
$c.Value = "7.46"
$c.ClearFormats()
$ws.Cells.Item(29,5).Value = "  -2.86%  "

$ws.Cells.Item(30,2).Value = "Binance-PegBSC-USD"
$ws.Cells.Item(30,3).Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Cells.Item(30,4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()
$ws.Cells.Item(30,5).Value = "  +0.27%  "

$c = $ws.Cells.Item(31,4)
$c.NumberFormat = "@"
$c.Value = "8.27"
$c.ClearFormats()
$ws.Cells.Item(31,5).Value = "  -0.47%  "

$c = $ws.Cells.Item(32,4)
$c.NumberFormat = "@"
$c.Value = "2.23"
$c.ClearFormats()
$ws.Cells.Item(32,5).Value = "  -2.59%  "

$ws.Cells.Item(33,4).Value = "3.486.82"
$ws.Cells.Item(33,5).Value = "  -1.03%  "

$ws.Cells.Item(34,5).Value = "  -0.03%  "

$c = $ws.Cells.Item(35,4)
$c.NumberFormat = "@"
$c.Value = "23.41"
$c.ClearFormats()
$ws.Cells.Item(35,5).Value = "  -2.61%  "

$ws.Cells.Item(36,5).Value = "  -0.38%  "

$ws.Cells.Item(37,5).Value = "  +2.83%  "

$ws.Cells.Item(38,5).Value = "  +0.20%  "

$c = $ws.Cells.Item(39,4)
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.ClearFormats()
$ws.Cells.Item(39,5).Value = "  +2.22%  "

$c = $ws.Cells.Item(40,4)
$c.NumberFormat = "@"
$c.Value = "161.74"
$c.ClearFormats()
$ws.Cells.Item(40,5).Value = "  -4.05%  "

$ws.Cells.Item(41,5).Value = "  -2.37%  "

$ws.Cells.Item(42,2).Value = "EnergySwap"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Cells.Item(42,4)
$c.NumberFormat = "@"
$c.Value = "26.53"
$c.ClearFormats()
$ws.Cells.Item(42,5).Value = "  +5.03%  "

$ws.Cells.Item(43,2).Value = "Mantle"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Cells.Item(43,4)
$c.NumberFormat = "@"
$c.Value = "0.808"
$c.ClearFormats()
$ws.Cells.Item(43,5).Value = "  -1.50%  "

$ws.Cells.Item(44,5).Value = "  -0.08%  "

$ws.Cells.Item(45,5).Value = "  -2.52%  "

$c = $ws.Cells.Item(46,4)
$c.NumberFormat = "@"
$c.Value = "41.37"
$c.ClearFormats()
$ws.Cells.Item(46,5).Value = "  -3.08%  "

$ws.Cells.Item(47,5).Value = "  -0.73%  "

$ws.Cells.Item(48,5).Value = "  -1.56%  "

$ws.Cells.Item(49,2).Value = "Maker"
$ws.Cells.Item(49,3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(49,4).Value = "2.418.91"
$ws.Cells.Item(49,5).Value = "  +0.56%  "

$ws.Cells.Item(50,2).Value = "Cosmos"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Cells.Item(50,4)
$c.NumberFormat = "@"
$c.Value = "6.82"
$c.ClearFormats()
$ws.Cells.Item(50,5).Value = "  -0.93%  "

$c = $ws.Cells.Item(51,4)
$c.NumberFormat = "@"
$c.Value = "0.891"
$c.ClearFormats()
$ws.Cells.Item(51,5).Value = "  -0.06%  "
